# "unify the conception of DataNode, DataTable, Entity."
#
# Rename the two sheets (formerly "Property1"/"Property2") to the unified
# "DataNode_*" naming scheme, nudge the header-row heights to their
# resaved values, and leave the workbook with the second sheet
# active/selected (matching the saved state captured in the target file).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "DataNode_1"
$ws2.Name = "DataNode_2"

# Header row (title row) and the type/comment row just above the data rows
# both come out a bit shorter/taller after the resave -- match the new
# row heights on both sheets.
$ws1.Rows.Item(1).RowHeight = 27
$ws1.Rows.Item(8).RowHeight = 40.5

$ws2.Rows.Item(1).RowHeight = 27
$ws2.Rows.Item(8).RowHeight = 67.5

# The workbook was left with the second sheet ("DataNode_2") as the active
# tab (previously the first sheet was active).
$ws2.Activate()
